$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.NumberFormat = "General"
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "30.522.32"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.883.65"
$ws.Range("E3").Value = "  +0.72%  "
Set-TextValue "D4" "1.000"
Set-TextValue "D5" "243.82"
$ws.Range("E5").Value = "  -1.55%  "
Set-TextValue "D6" "1.000"
$ws.Range("E6").Value = "  +0.03%  "
Set-TextValue "D7" "0.4684"
$ws.Range("E7").Value = "  -1.08%  "
Set-TextValue "D8" "0.2893"
$ws.Range("E8").Value = "  -0.61%  "
Set-TextValue "D9" "0.06482"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +1.21%  "
Set-TextValue "D11" "0.07741"
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "1.883.48"
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D13" "0.7293"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("B14").Value = "Litecoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D14" "95.37"
$ws.Range("E14").Value = "  -1.15%  "
Set-TextValue "D15" "5.181"
$ws.Range("E15").Value = "  +0.83%  "
Set-TextValue "D16" "282.46"
$ws.Range("E16").Value = "  +3.88%  "
$ws.Range("D17").Value = "30.513.16"
$ws.Range("E17").Value = "  -0.18%  "
Set-TextValue "D18" "13.01"
$ws.Range("E18").Value = "  -2.02%  "
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  +0.01%  "
Set-TextValue "D20" "0.000007474"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "2.128.33"
$ws.Range("E21").Value = "  +0.68%  "
Set-TextValue "D22" "0.9999"
$ws.Range("E22").Value = "  -0.03%  "
Set-TextValue "D23" "5.252"
$ws.Range("E23").Value = "  +0.29%  "
Set-TextValue "D24" "6.252"
$ws.Range("E24").Value = "  +1.51%  "
Set-TextValue "D25" "163.62"
$ws.Range("E25").Value = "  -0.16%  "
Set-TextValue "D26" "9.081"
$ws.Range("E26").Value = "  -1.22%  "
Set-TextValue "D27" "18.88"
$ws.Range("E27").Value = "  +0.74%  "
Set-TextValue "D28" "1.891"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  -0.96%  "
Set-TextValue "D30" "0.09705"
$ws.Range("E30").Value = "  -2.65%  "
Set-TextValue "D31" "1.472"
$ws.Range("E31").Value = "  -2.47%  "
Set-TextValue "D32" "4.275"
$ws.Range("E32").Value = "  +0.09%  "
Set-TextValue "D33" "4.123"
$ws.Range("E33").Value = "  +0.67%  "
Set-TextValue "D34" "0.04856"
$ws.Range("E34").Value = "  +1.69%  "
Set-TextValue "D35" "1.124"
$ws.Range("E35").Value = "  +0.60%  "
Set-TextValue "D36" "0.6922"
$ws.Range("E36").Value = "  -0.26%  "
Set-TextValue "D37" "2.722"
$ws.Range("E37").Value = "  +0.17%  "
Set-TextValue "D38" "0.01890"
$ws.Range("E38").Value = "  +2.32%  "
Set-TextValue "D39" "2.823"
$ws.Range("E39").Value = "  +2.76%  "
Set-TextValue "D40" "75.52"
$ws.Range("E40").Value = "  +3.74%  "
Set-TextValue "D41" "6.167"
$ws.Range("E41").Value = "  -0.02%  "
Set-TextValue "D42" "2.011"
$ws.Range("E42").Value = "  +2.45%  "
Set-TextValue "D43" "0.4250"
$ws.Range("E43").Value = "  +2.02%  "
$ws.Range("E44").Value = "  -0.01%  "
Set-TextValue "D45" "0.8226"
$ws.Range("E45").Value = "  -1.12%  "
Set-TextValue "D46" "101.28"
$ws.Range("E46").Value = "  -0.22%  "
Set-TextValue "D47" "9.508"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D48" "6.962"
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D49" "35.21"
$ws.Range("E49").Value = "  -0.59%  "
Set-TextValue "D50" "914.72"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("E51").Value = "  +1.90%  "
